# "custom accuracy + 데이터 1000개"
#
# 1) Re-express every numeric value in row 5 (columns B..AH) at "custom
#    accuracy" (2 decimal places) -- this is not a uniform ROUND() of the
#    in-sheet figures (their source values round inconsistently around
#    half-way points), so each column is set to its known, already-rounded
#    figure directly.
# 2) Drop the old row 6, trimming the sample down and pulling the sheet's
#    used range back in from A1:AH6 to A1:AH5.
# 3) Column AG (the 33rd column) narrows from width 8 to width 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5Values = [ordered]@{
    "B5"  = 24.56
    "C5"  = 18.69
    "D5"  = 0.7
    "E5"  = 54.01
    "F5"  = 44.55
    "G5"  = 19.54
    "H5"  = 75.69
    "I5"  = 30.31
    "J5"  = 13.74
    "K5"  = 20.58
    "L5"  = 22.03
    "M5"  = 22.99
    "N5"  = 6.28
    "O5"  = 19.33
    "P5"  = 27.89
    "Q5"  = 16.21
    "R5"  = 0.27
    "S5"  = 0.81
    "T5"  = 290.11
    "U5"  = 54.67
    "V5"  = 17.94
    "W5"  = 36.73
    "X5"  = 19.75
    "Y5"  = 2.56
    "Z5"  = 37.04
    "AA5" = 15.92
    "AB5" = 14.4
    "AC5" = 16.94
    "AD5" = 23.16
    "AE5" = 0.21
    "AF5" = 68.71
    "AG5" = 10.48
    "AH5" = 22.38
}

foreach ($addr in $row5Values.Keys) {
    $ws.Range($addr).Value = $row5Values[$addr]
}

# Remove the now-unneeded last data row (old row 6); the sheet's dimension
# collapses from A1:AH6 to A1:AH5 as a result.
$ws.Rows(6).Delete()

# Narrow column AG (the 33rd column) from raw width 8 to raw width 7.
# Excel's ColumnWidth COM property is offset from the raw OOXML <col width>
# value by a constant padding factor, so shave one unit off the current
# COM-reported width rather than hard-coding the converted figure.
$ws.Columns(33).ColumnWidth = $ws.Columns(33).ColumnWidth - 1
